$wb = $excel.ActiveWorkbook

# --- Part 1: convert D517:D525 on the "day" sheet from text to numeric ---
$wsDay = $wb.Worksheets.Item("day")
$wsDay.Cells.Item(517, 4).Value = 500027
$wsDay.Cells.Item(518, 4).Value = 500790
$wsDay.Cells.Item(519, 4).Value = 506395
$wsDay.Cells.Item(520, 4).Value = 500085
$wsDay.Cells.Item(521, 4).Value = 500086
$wsDay.Cells.Item(522, 4).Value = 532134
$wsDay.Cells.Item(523, 4).Value = 540065
$wsDay.Cells.Item(524, 4).Value = 541153
$wsDay.Cells.Item(525, 4).Value = 532461

# --- Part 2: append 30 new rows (260:289) to the "week" sheet ---
$wsWeek = $wb.Worksheets.Item("week")

# Column D on the new rows must stay text (stock-code strings), like the existing data.
# Force text storage via NumberFormat, then restore the default style so no stray
# number-format leaks onto the cells themselves.
$dRange = $wsWeek.Range("D260:D289")
$dRange.NumberFormat = "@"

# row 260
$wsWeek.Cells.Item(260, 1).Value = 1
$wsWeek.Cells.Item(260, 2).Value = "MRF"
$wsWeek.Cells.Item(260, 3).Value = "Mrf Limited"
$wsWeek.Cells.Item(260, 4).Value = "500290"
$wsWeek.Cells.Item(260, 5).Value = -1.29
$wsWeek.Cells.Item(260, 6).Value = 134260.75
$wsWeek.Cells.Item(260, 7).Value = 3820
$wsWeek.Cells.Item(260, 8).Value = "week"
$wsWeek.Cells.Item(260, 9).Value = "06/09/2024 11:32:28"

# row 261
$wsWeek.Cells.Item(261, 1).Value = 2
$wsWeek.Cells.Item(261, 2).Value = "PAGEIND"
$wsWeek.Cells.Item(261, 3).Value = "Page Industries Limited"
$wsWeek.Cells.Item(261, 4).Value = "532827"
$wsWeek.Cells.Item(261, 5).Value = -1.69
$wsWeek.Cells.Item(261, 6).Value = 40360.4
$wsWeek.Cells.Item(261, 7).Value = 17008
$wsWeek.Cells.Item(261, 8).Value = "week"
$wsWeek.Cells.Item(261, 9).Value = "06/09/2024 11:32:28"

# row 262
$wsWeek.Cells.Item(262, 1).Value = 3
$wsWeek.Cells.Item(262, 2).Value = "MARUTI"
$wsWeek.Cells.Item(262, 3).Value = "Maruti Suzuki India Limited"
$wsWeek.Cells.Item(262, 4).Value = "532500"
$wsWeek.Cells.Item(262, 5).Value = -0.91
$wsWeek.Cells.Item(262, 6).Value = 12186.15
$wsWeek.Cells.Item(262, 7).Value = 733354
$wsWeek.Cells.Item(262, 8).Value = "week"
$wsWeek.Cells.Item(262, 9).Value = "06/09/2024 11:32:28"

# row 263
$wsWeek.Cells.Item(263, 1).Value = 4
$wsWeek.Cells.Item(263, 2).Value = "LT"
$wsWeek.Cells.Item(263, 3).Value = "Larsen & Toubro Limited"
$wsWeek.Cells.Item(263, 4).Value = "500510"
$wsWeek.Cells.Item(263, 5).Value = -1.36
$wsWeek.Cells.Item(263, 6).Value = 3574.75
$wsWeek.Cells.Item(263, 7).Value = 3096557
$wsWeek.Cells.Item(263, 8).Value = "week"
$wsWeek.Cells.Item(263, 9).Value = "06/09/2024 11:32:28"

# row 264
$wsWeek.Cells.Item(264, 1).Value = 5
$wsWeek.Cells.Item(264, 2).Value = "NAVINFLUOR"
$wsWeek.Cells.Item(264, 3).Value = "Navin Fluorine International Limited"
$wsWeek.Cells.Item(264, 4).Value = "532504"
$wsWeek.Cells.Item(264, 5).Value = -1.49
$wsWeek.Cells.Item(264, 6).Value = 3298.35
$wsWeek.Cells.Item(264, 7).Value = 174121
$wsWeek.Cells.Item(264, 8).Value = "week"
$wsWeek.Cells.Item(264, 9).Value = "06/09/2024 11:32:28"

# row 265
$wsWeek.Cells.Item(265, 1).Value = 6
$wsWeek.Cells.Item(265, 2).Value = "INDIAMART"
$wsWeek.Cells.Item(265, 3).Value = "Indiamart Intermesh Ltd"
$wsWeek.Cells.Item(265, 4).Value = "542726"
$wsWeek.Cells.Item(265, 5).Value = -1.5
$wsWeek.Cells.Item(265, 6).Value = 2983.8
$wsWeek.Cells.Item(265, 7).Value = 167359
$wsWeek.Cells.Item(265, 8).Value = "week"
$wsWeek.Cells.Item(265, 9).Value = "06/09/2024 11:32:28"

# row 266
$wsWeek.Cells.Item(266, 1).Value = 7
$wsWeek.Cells.Item(266, 2).Value = "ADANIENT"
$wsWeek.Cells.Item(266, 3).Value = "Adani Enterprises Limited"
$wsWeek.Cells.Item(266, 4).Value = "512599"
$wsWeek.Cells.Item(266, 5).Value = -1.32
$wsWeek.Cells.Item(266, 6).Value = 2975.45
$wsWeek.Cells.Item(266, 7).Value = 1251165
$wsWeek.Cells.Item(266, 8).Value = "week"
$wsWeek.Cells.Item(266, 9).Value = "06/09/2024 11:32:28"

# row 267
$wsWeek.Cells.Item(267, 1).Value = 8
$wsWeek.Cells.Item(267, 2).Value = "BALKRISIND"
$wsWeek.Cells.Item(267, 3).Value = "Balkrishna Industries Limited"
$wsWeek.Cells.Item(267, 4).Value = "502355"
$wsWeek.Cells.Item(267, 5).Value = 0.55
$wsWeek.Cells.Item(267, 6).Value = 2973.7
$wsWeek.Cells.Item(267, 7).Value = 282060
$wsWeek.Cells.Item(267, 8).Value = "week"
$wsWeek.Cells.Item(267, 9).Value = "06/09/2024 11:32:28"

# row 268
$wsWeek.Cells.Item(268, 1).Value = 9
$wsWeek.Cells.Item(268, 2).Value = "DEEPAKNTR"
$wsWeek.Cells.Item(268, 3).Value = "Deepak Nitrite Limited"
$wsWeek.Cells.Item(268, 4).Value = "506401"
$wsWeek.Cells.Item(268, 5).Value = -2.32
$wsWeek.Cells.Item(268, 6).Value = 2929.1
$wsWeek.Cells.Item(268, 7).Value = 220419
$wsWeek.Cells.Item(268, 8).Value = "week"
$wsWeek.Cells.Item(268, 9).Value = "06/09/2024 11:32:28"

# row 269
$wsWeek.Cells.Item(269, 1).Value = 10
$wsWeek.Cells.Item(269, 2).Value = "M&M"
$wsWeek.Cells.Item(269, 3).Value = "Mahindra & Mahindra Limited"
$wsWeek.Cells.Item(269, 4).Value = "500520"
$wsWeek.Cells.Item(269, 5).Value = -0.92
$wsWeek.Cells.Item(269, 6).Value = 2698.1
$wsWeek.Cells.Item(269, 7).Value = 2835519
$wsWeek.Cells.Item(269, 8).Value = "week"
$wsWeek.Cells.Item(269, 9).Value = "06/09/2024 11:32:28"

# row 270
$wsWeek.Cells.Item(270, 1).Value = 11
$wsWeek.Cells.Item(270, 2).Value = "SRF"
$wsWeek.Cells.Item(270, 3).Value = "Srf Limited"
$wsWeek.Cells.Item(270, 4).Value = "503806"
$wsWeek.Cells.Item(270, 5).Value = -4.18
$wsWeek.Cells.Item(270, 6).Value = 2509.05
$wsWeek.Cells.Item(270, 7).Value = 714722
$wsWeek.Cells.Item(270, 8).Value = "week"
$wsWeek.Cells.Item(270, 9).Value = "06/09/2024 11:32:28"

# row 271
$wsWeek.Cells.Item(271, 1).Value = 12
$wsWeek.Cells.Item(271, 2).Value = "NESTLEIND"
$wsWeek.Cells.Item(271, 3).Value = "Nestle India Limited"
$wsWeek.Cells.Item(271, 4).Value = "500790"
$wsWeek.Cells.Item(271, 5).Value = -0.07000000000000001
$wsWeek.Cells.Item(271, 6).Value = 2503.2
$wsWeek.Cells.Item(271, 7).Value = 836826
$wsWeek.Cells.Item(271, 8).Value = "week"
$wsWeek.Cells.Item(271, 9).Value = "06/09/2024 11:32:28"

# row 272
$wsWeek.Cells.Item(272, 1).Value = 13
$wsWeek.Cells.Item(272, 2).Value = "MGL"
$wsWeek.Cells.Item(272, 3).Value = "Mahanagar Gas Limited"
$wsWeek.Cells.Item(272, 4).Value = "539957"
$wsWeek.Cells.Item(272, 5).Value = -2.01
$wsWeek.Cells.Item(272, 6).Value = 1839
$wsWeek.Cells.Item(272, 7).Value = 424620
$wsWeek.Cells.Item(272, 8).Value = "week"
$wsWeek.Cells.Item(272, 9).Value = "06/09/2024 11:32:28"

# row 273
$wsWeek.Cells.Item(273, 1).Value = 14
$wsWeek.Cells.Item(273, 2).Value = "ADANIPORTS"
$wsWeek.Cells.Item(273, 3).Value = "Adani Ports And Special Economic Zone Limited"
$wsWeek.Cells.Item(273, 4).Value = "532921"
$wsWeek.Cells.Item(273, 5).Value = -1.58
$wsWeek.Cells.Item(273, 6).Value = 1442.4
$wsWeek.Cells.Item(273, 7).Value = 3340767
$wsWeek.Cells.Item(273, 8).Value = "week"
$wsWeek.Cells.Item(273, 9).Value = "06/09/2024 11:32:28"

# row 274
$wsWeek.Cells.Item(274, 1).Value = 15
$wsWeek.Cells.Item(274, 2).Value = "BATAINDIA"
$wsWeek.Cells.Item(274, 3).Value = "Bata India Limited"
$wsWeek.Cells.Item(274, 4).Value = "500043"
$wsWeek.Cells.Item(274, 5).Value = -2.13
$wsWeek.Cells.Item(274, 6).Value = 1406.25
$wsWeek.Cells.Item(274, 7).Value = 395504
$wsWeek.Cells.Item(274, 8).Value = "week"
$wsWeek.Cells.Item(274, 9).Value = "06/09/2024 11:32:28"

# row 275
$wsWeek.Cells.Item(275, 1).Value = 16
$wsWeek.Cells.Item(275, 2).Value = "TATACHEM"
$wsWeek.Cells.Item(275, 3).Value = "Tata Chemicals Limited"
$wsWeek.Cells.Item(275, 4).Value = "500770"
$wsWeek.Cells.Item(275, 5).Value = -2.28
$wsWeek.Cells.Item(275, 6).Value = 1056.65
$wsWeek.Cells.Item(275, 7).Value = 1014325
$wsWeek.Cells.Item(275, 8).Value = "week"
$wsWeek.Cells.Item(275, 9).Value = "06/09/2024 11:32:28"

# row 276
$wsWeek.Cells.Item(276, 1).Value = 17
$wsWeek.Cells.Item(276, 2).Value = "TATAMOTORS"
$wsWeek.Cells.Item(276, 3).Value = "Tata Motors Limited"
$wsWeek.Cells.Item(276, 4).Value = "500570"
$wsWeek.Cells.Item(276, 5).Value = -1.85
$wsWeek.Cells.Item(276, 6).Value = 1049.35
$wsWeek.Cells.Item(276, 7).Value = 8717377
$wsWeek.Cells.Item(276, 8).Value = "week"
$wsWeek.Cells.Item(276, 9).Value = "06/09/2024 11:32:28"

# row 277
$wsWeek.Cells.Item(277, 1).Value = 18
$wsWeek.Cells.Item(277, 2).Value = "DLF"
$wsWeek.Cells.Item(277, 3).Value = "Dlf Limited"
$wsWeek.Cells.Item(277, 4).Value = "532868"
$wsWeek.Cells.Item(277, 5).Value = -3.26
$wsWeek.Cells.Item(277, 6).Value = 814.25
$wsWeek.Cells.Item(277, 7).Value = 3788193
$wsWeek.Cells.Item(277, 8).Value = "week"
$wsWeek.Cells.Item(277, 9).Value = "06/09/2024 11:32:28"

# row 278
$wsWeek.Cells.Item(278, 1).Value = 19
$wsWeek.Cells.Item(278, 2).Value = "APOLLOTYRE"
$wsWeek.Cells.Item(278, 3).Value = "Apollo Tyres Limited"
$wsWeek.Cells.Item(278, 4).Value = "500877"
$wsWeek.Cells.Item(278, 5).Value = -0.83
$wsWeek.Cells.Item(278, 6).Value = 507.75
$wsWeek.Cells.Item(278, 7).Value = 2355690
$wsWeek.Cells.Item(278, 8).Value = "week"
$wsWeek.Cells.Item(278, 9).Value = "06/09/2024 11:32:28"

# row 279
$wsWeek.Cells.Item(279, 1).Value = 20
$wsWeek.Cells.Item(279, 2).Value = "TATAPOWER"
$wsWeek.Cells.Item(279, 3).Value = "Tata Power Company Limited"
$wsWeek.Cells.Item(279, 4).Value = "500400"
$wsWeek.Cells.Item(279, 5).Value = -0.93
$wsWeek.Cells.Item(279, 6).Value = 417
$wsWeek.Cells.Item(279, 7).Value = 10236248
$wsWeek.Cells.Item(279, 8).Value = "week"
$wsWeek.Cells.Item(279, 9).Value = "06/09/2024 11:32:28"

# row 280
$wsWeek.Cells.Item(280, 1).Value = 21
$wsWeek.Cells.Item(280, 2).Value = "NTPC"
$wsWeek.Cells.Item(280, 3).Value = "Ntpc Limited"
$wsWeek.Cells.Item(280, 4).Value = "532555"
$wsWeek.Cells.Item(280, 5).Value = -2.1
$wsWeek.Cells.Item(280, 6).Value = 394.8
$wsWeek.Cells.Item(280, 7).Value = 17956132
$wsWeek.Cells.Item(280, 8).Value = "week"
$wsWeek.Cells.Item(280, 9).Value = "06/09/2024 11:32:28"

# row 281
$wsWeek.Cells.Item(281, 1).Value = 22
$wsWeek.Cells.Item(281, 2).Value = "INDIACEM"
$wsWeek.Cells.Item(281, 3).Value = "The India Cements Limited"
$wsWeek.Cells.Item(281, 4).Value = "530005"
$wsWeek.Cells.Item(281, 5).Value = -0.3
$wsWeek.Cells.Item(281, 6).Value = 363.35
$wsWeek.Cells.Item(281, 7).Value = 1207203
$wsWeek.Cells.Item(281, 8).Value = "week"
$wsWeek.Cells.Item(281, 9).Value = "06/09/2024 11:32:28"

# row 282
$wsWeek.Cells.Item(282, 1).Value = 23
$wsWeek.Cells.Item(282, 2).Value = "POWERGRID"
$wsWeek.Cells.Item(282, 3).Value = "Power Grid Corporation Of India Limited"
$wsWeek.Cells.Item(282, 4).Value = "532898"
$wsWeek.Cells.Item(282, 5).Value = -0.44
$wsWeek.Cells.Item(282, 6).Value = 329.8
$wsWeek.Cells.Item(282, 7).Value = 23755639
$wsWeek.Cells.Item(282, 8).Value = "week"
$wsWeek.Cells.Item(282, 9).Value = "06/09/2024 11:32:28"

# row 283
$wsWeek.Cells.Item(283, 1).Value = 24
$wsWeek.Cells.Item(283, 2).Value = "ABFRL"
$wsWeek.Cells.Item(283, 3).Value = "Aditya Birla Fashion And Retail Limited"
$wsWeek.Cells.Item(283, 4).Value = "535755"
$wsWeek.Cells.Item(283, 5).Value = -1.95
$wsWeek.Cells.Item(283, 6).Value = 309.15
$wsWeek.Cells.Item(283, 7).Value = 2426890
$wsWeek.Cells.Item(283, 8).Value = "week"
$wsWeek.Cells.Item(283, 9).Value = "06/09/2024 11:32:28"

# row 284
$wsWeek.Cells.Item(284, 1).Value = 25
$wsWeek.Cells.Item(284, 2).Value = "GAIL"
$wsWeek.Cells.Item(284, 3).Value = "Gail (india) Limited"
$wsWeek.Cells.Item(284, 4).Value = "532155"
$wsWeek.Cells.Item(284, 5).Value = -2.32
$wsWeek.Cells.Item(284, 6).Value = 222.82
$wsWeek.Cells.Item(284, 7).Value = 12081433
$wsWeek.Cells.Item(284, 8).Value = "week"
$wsWeek.Cells.Item(284, 9).Value = "06/09/2024 11:32:28"

# row 285
$wsWeek.Cells.Item(285, 1).Value = 26
$wsWeek.Cells.Item(285, 2).Value = "BANDHANBNK"
$wsWeek.Cells.Item(285, 3).Value = "Bandhan Bank Ltd"
$wsWeek.Cells.Item(285, 4).Value = "541153"
$wsWeek.Cells.Item(285, 5).Value = -3.6
$wsWeek.Cells.Item(285, 6).Value = 196.33
$wsWeek.Cells.Item(285, 7).Value = 19496566
$wsWeek.Cells.Item(285, 8).Value = "week"
$wsWeek.Cells.Item(285, 9).Value = "06/09/2024 11:32:28"

# row 286
$wsWeek.Cells.Item(286, 1).Value = 27
$wsWeek.Cells.Item(286, 2).Value = "IOC"
$wsWeek.Cells.Item(286, 3).Value = "Indian Oil Corporation Limited"
$wsWeek.Cells.Item(286, 4).Value = "530965"
$wsWeek.Cells.Item(286, 5).Value = -2.59
$wsWeek.Cells.Item(286, 6).Value = 176.64
$wsWeek.Cells.Item(286, 7).Value = 31879056
$wsWeek.Cells.Item(286, 8).Value = "week"
$wsWeek.Cells.Item(286, 9).Value = "06/09/2024 11:32:28"

# row 287
$wsWeek.Cells.Item(287, 1).Value = 28
$wsWeek.Cells.Item(287, 2).Value = "CUB"
$wsWeek.Cells.Item(287, 3).Value = "City Union Bank Limited"
$wsWeek.Cells.Item(287, 4).Value = "532210"
$wsWeek.Cells.Item(287, 5).Value = -2.72
$wsWeek.Cells.Item(287, 6).Value = 166.42
$wsWeek.Cells.Item(287, 7).Value = 3914747
$wsWeek.Cells.Item(287, 8).Value = "week"
$wsWeek.Cells.Item(287, 9).Value = "06/09/2024 11:32:28"

# row 288
$wsWeek.Cells.Item(288, 1).Value = 29
$wsWeek.Cells.Item(288, 2).Value = "PNB"
$wsWeek.Cells.Item(288, 3).Value = "Punjab National Bank"
$wsWeek.Cells.Item(288, 4).Value = "532461"
$wsWeek.Cells.Item(288, 5).Value = -3
$wsWeek.Cells.Item(288, 6).Value = 110
$wsWeek.Cells.Item(288, 7).Value = 29092655
$wsWeek.Cells.Item(288, 8).Value = "week"
$wsWeek.Cells.Item(288, 9).Value = "06/09/2024 11:32:28"

# row 289
$wsWeek.Cells.Item(289, 1).Value = 30
$wsWeek.Cells.Item(289, 2).Value = "GMRINFRA"
$wsWeek.Cells.Item(289, 3).Value = "Gmr Infrastructure Limited"
$wsWeek.Cells.Item(289, 4).Value = "532754"
$wsWeek.Cells.Item(289, 5).Value = -4.95
$wsWeek.Cells.Item(289, 6).Value = 91.03
$wsWeek.Cells.Item(289, 7).Value = 39150062
$wsWeek.Cells.Item(289, 8).Value = "week"
$wsWeek.Cells.Item(289, 9).Value = "06/09/2024 11:32:28"

# Reset the style on the D column back to Normal/default so no extra number format
# sticks to these cells (keeps styles.xml untouched, matching the source data).
$dRange.Style = "Normal"
